# Fruta / hortaliza, semanal
#
# A new week's worth of price reports (3 rows, dated 2021-11-11 / serial 44511)
# is inserted right above the existing "Hayward" block that used to start at
# row 396. That push shifts the old rows 396:490 down to 399:493 (identical
# content - no re-typing needed, Excel's own Insert does that), and the sheet
# dimension grows from A1:T490 to A1:T493.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 396

# Insert 3 blank rows above the current row 396, pushing 396:490 -> 399:493.
$ws.Rows.Item($startRow).Resize(3).Insert()

# The 3 freshly inserted rows all share the same "static" columns (A-K, Q, R, T)
# as the rows that used to sit there (now 399:401); only the date (D), grade
# (L), price (M/N/O/P) and per-unit price (S) differ per grade.
$newRows = @(
    @{ L = "Especial"; M = 60; N = 13000; S = 1300 },
    @{ L = "Primera";  M = 65; N = 12000; S = 1200 },
    @{ L = "Segunda";  M = 60; N = 10000; S = 1000 }
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($r, 1).Value = 3                              # A: Mercado ID
    $ws.Cells.Item($r, 2).Value = "Femacal de La Calera"          # B: Mercado
    $ws.Cells.Item($r, 3).Value = "Coquimbo"                      # C: Región Mercado
    $ws.Cells.Item($r, 4).Value = 44511                           # D: Fecha (2021-11-11)
    $ws.Cells.Item($r, 5).Value = 5                               # E
    $ws.Cells.Item($r, 6).Value = "Fruta"                         # F
    $ws.Cells.Item($r, 7).Value = 100101                          # G
    $ws.Cells.Item($r, 8).Value = "Berries"                       # H
    $ws.Cells.Item($r, 9).Value = 100101007                       # I
    $ws.Cells.Item($r, 10).Value = "Kiwi"                         # J
    $ws.Cells.Item($r, 11).Value = "Hayward"                      # K: Variedad
    $ws.Cells.Item($r, 12).Value = $data.L                        # L: Calidad
    $ws.Cells.Item($r, 13).Value = $data.M                        # M: Cantidad
    $ws.Cells.Item($r, 14).Value = $data.N                        # N: Precio minimo
    $ws.Cells.Item($r, 15).Value = $data.N                        # O: Precio corriente
    $ws.Cells.Item($r, 16).Value = $data.N                        # P: Precio maximo
    $ws.Cells.Item($r, 17).Value = "$/bandeja 10 kilos"           # Q: Unidad
    $ws.Cells.Item($r, 18).Value = "Región de O'Higgins"          # R: Origen
    $ws.Cells.Item($r, 19).Value = $data.S                        # S: Precio por kilo
    $ws.Cells.Item($r, 20).Value = 10                             # T: Kilos por unidad
}
